$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.836640596389771
$ws.Range("B1").Value = 4.23958158493042
$ws.Range("C1").Value = 2.800498962402344
$ws.Range("D1").Value = 2.319116592407227
$ws.Range("E1").Value = 1.872685670852661
